# Fixed issues with 81RF protective element: change the default 81x
# (frequency) protection settings so they are effectively disabled for
# every relay row on the "relays" sheet.
#
# Columns T:V hold the primary-side 81x pickup/time/definite-time settings
# and columns AH:AJ hold the matching backup-side settings. Both blocks get
# the same new defaults: pickup 100, time dial 10, definite time 0.1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

for ($row = 2; $row -le 7; $row++) {
    $ws.Range("T" + $row).Value = 100
    $ws.Range("U" + $row).Value = 10
    $ws.Range("V" + $row).Value = 0.1

    $ws.Range("AH" + $row).Value = 100
    $ws.Range("AI" + $row).Value = 10
    $ws.Range("AJ" + $row).Value = 0.1
}

# Leave the selection where the author last left it while reviewing the
# backup 81x block.
$ws.Activate()
$ws.Range("AH2:AJ7").Select()
